$wb = $excel.ActiveWorkbook

# Map of worksheet name -> (DateProd cell B2, DateDemo cell D2) values
# These reflect a second Katalon bootstrap test run being recorded:
# column B = DateProd, column D = DateDemo timestamps.

$updates = @{
    "AddDeleteRole"  = @{ B2 = "Mon Nov 10 16:09:18 IST 2025"; D2 = "Thu Nov 06 22:15:07 IST 2025" }
    "SearchRole"     = @{ B2 = "Mon Nov 10 16:09:58 IST 2025"; D2 = "Thu Nov 06 22:15:50 IST 2025" }
    "CreateUser"     = @{ B2 = "Mon Nov 10 16:10:32 IST 2025"; D2 = "Thu Nov 06 22:16:26 IST 2025" }
    "FindUser"       = @{ B2 = "Mon Nov 10 16:11:08 IST 2025"; D2 = "Thu Nov 06 22:17:05 IST 2025" }
    "ModifyUser"     = @{ B2 = "Mon Nov 10 16:11:42 IST 2025"; D2 = "Thu Nov 06 22:17:42 IST 2025" }
    "ModifyUserPwd"  = @{ B2 = "Mon Nov 10 16:12:33 IST 2025"; D2 = "Thu Nov 06 22:18:29 IST 2025" }
    "FindCaseUser"   = @{ B2 = "Mon Nov 10 16:13:47 IST 2025"; D2 = "Thu Nov 06 22:19:50 IST 2025" }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $updates[$sheetName]
    $ws.Range("B2").Value = $vals.B2
    $ws.Range("D2").Value = $vals.D2
}
